# "added some unit tests" - extends the importer test fixture sheet with
# one more data column (I), giving every existing row (1-7) a constant
# 6.66 value in that new column, and leaves the selection on the new
# last cell (I7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1:I7").Value = 6.66

$ws.Range("I7").Select()
